$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.72694269319203
$ws.Range("C2").Value = 10.02214985923283
$ws.Range("E2").Value = 15.67841776750739
$ws.Range("F2").Value = 40.29838067125225
$ws.Range("G2").Value = 3.674794739095985
$ws.Range("J2").Value = 8.481277689131128
$ws.Range("K2").Value = 10.30328710115067
$ws.Range("L2").Value = 12.19591756634791
$ws.Range("N2").Value = 20.43505483960072
$ws.Range("O2").Value = 26.8640291203732

$ws.Range("B3").Value = 14.52577575470528
$ws.Range("C3").Value = 10.03588197905311
$ws.Range("E3").Value = 15.6749958021157
$ws.Range("F3").Value = 40.31766005521795
$ws.Range("G3").Value = 3.676536782166882
$ws.Range("J3").Value = 8.480015979050343
$ws.Range("K3").Value = 10.16472380428016
$ws.Range("L3").Value = 12.17689187081391
$ws.Range("N3").Value = 20.49608661225827
$ws.Range("O3").Value = 26.94583676861788

$ws.Range("B4").Value = 14.40361127369492
$ws.Range("C4").Value = 10.04481949068168
$ws.Range("E4").Value = 15.67528297528058
$ws.Range("F4").Value = 40.33809977951206
$ws.Range("G4").Value = 3.677663511670437
$ws.Range("J4").Value = 8.479325010176394
$ws.Range("K4").Value = 10.08037332212831
$ws.Range("L4").Value = 12.16694818076342
$ws.Range("N4").Value = 20.53531187402903
$ws.Range("O4").Value = 27.00084713093472

$ws.Range("B5").Value = 14.35422767855418
$ws.Range("C5").Value = 10.04858917193937
$ws.Range("E5").Value = 15.67600239067981
$ws.Range("F5").Value = 40.34859270290187
$ws.Range("G5").Value = 3.678137067592289
$ws.Range("J5").Value = 8.479064603694923
$ws.Range("K5").Value = 10.046222155222
$ws.Range("L5").Value = 12.16333629386287
$ws.Range("N5").Value = 20.55173822171144
$ws.Range("O5").Value = 27.02446510894241

$ws.Range("B6").Value = 14.34605341642152
$ws.Range("C6").Value = 10.04922284007987
$ws.Range("E6").Value = 15.6761582810294
$ws.Range("F6").Value = 40.35046571496124
$ws.Range("G6").Value = 3.678216572490033
$ws.Range("J6").Value = 8.479022645305175
$ws.Range("K6").Value = 10.04056595270685
$ws.Range("L6").Value = 12.16276322013059
$ws.Range("N6").Value = 20.55449252232216
$ws.Range("O6").Value = 27.02845934635881

$ws.Range("B7").Value = 14.4029435740891
$ws.Range("C7").Value = 10.04486981291329
$ws.Range("E7").Value = 15.67529023625487
$ws.Range("F7").Value = 40.33823253059621
$ws.Range("G7").Value = 3.67766983982892
$ws.Range("J7").Value = 8.479321412380783
$ws.Range("K7").Value = 10.0799117955555
$ws.Range("L7").Value = 12.16689768305254
$ws.Range("N7").Value = 20.53553161520357
$ws.Range("O7").Value = 27.00116079083966

$ws.Range("B8").Value = 14.65733520379532
$ws.Range("C8").Value = 10.0267798987001
$ws.Range("E8").Value = 15.6767433575893
$ws.Range("F8").Value = 40.30324341115117
$ws.Range("G8").Value = 3.675383569444044
$ws.Range("J8").Value = 8.480825342942488
$ws.Range("K8").Value = 10.25538300992692
$ws.Range("L8").Value = 12.18899872950491
$ws.Range("N8").Value = 20.45573585847206
$ws.Range("O8").Value = 26.89124357555349

$ws.Range("B9").Value = 15.16431819284102
$ws.Range("C9").Value = 9.995304144736474
$ws.Range("E9").Value = 15.69845377820867
$ws.Range("F9").Value = 40.30282545638072
$ws.Range("G9").Value = 3.671351307585515
$ws.Range("J9").Value = 8.48443405291898
$ws.Range("K9").Value = 10.60352921404752
$ws.Range("L9").Value = 12.24598098593008
$ws.Range("N9").Value = 20.31309348299676
$ws.Range("O9").Value = 26.71367674935442

$ws.Range("B10").Value = 15.53826905199952
$ws.Range("C10").Value = 9.974594580628562
$ws.Range("E10").Value = 15.72577092924718
$ws.Range("F10").Value = 40.34395993506232
$ws.Range("G10").Value = 3.668660984412944
$ws.Range("J10").Value = 8.487481299559759
$ws.Range("K10").Value = 10.85948189390988
$ws.Range("L10").Value = 12.29594777272384
$ws.Range("N10").Value = 20.21664321151805
$ws.Range("O10").Value = 26.60644198800437

$ws.Range("B11").Value = 15.70798540993317
$ws.Range("C11").Value = 9.965693176104121
$ws.Range("E11").Value = 15.74063330312708
$ws.Range("F11").Value = 40.37162242040339
$ws.Range("G11").Value = 3.667495586835761
$ws.Range("J11").Value = 8.48895223804997
$ws.Range("K11").Value = 10.9754819173111
$ws.Range("L11").Value = 12.32038714647999
$ws.Range("N11").Value = 20.17455999949638
$ws.Range("O11").Value = 26.56271494694356

$ws.Range("B12").Value = 15.7721370671525
$ws.Range("C12").Value = 9.962396794780865
$ws.Range("E12").Value = 15.74660829003183
$ws.Range("F12").Value = 40.3833784580677
$ws.Range("G12").Value = 3.667062639814848
$ws.Range("J12").Value = 8.489521307671485
$ws.Range("K12").Value = 11.01930708227928
$ws.Range("L12").Value = 12.32988272116686
$ws.Range("N12").Value = 20.15888058296382
$ws.Range("O12").Value = 26.54688460306303

$ws.Range("B13").Value = 15.75832695554898
$ws.Range("C13").Value = 9.963103426029994
$ws.Range("E13").Value = 15.74530609546986
$ws.Range("F13").Value = 40.38078972850429
$ws.Range("G13").Value = 3.667155511362383
$ws.Range("J13").Value = 8.489398214140669
$ws.Range("K13").Value = 11.00987367410983
$ws.Range("L13").Value = 12.32782704471202
$ws.Range("N13").Value = 20.16224603289982
$ws.Range("O13").Value = 26.55026154684615

$ws.Range("B14").Value = 15.71326586858218
$ws.Range("C14").Value = 9.96542049177607
$ws.Range("E14").Value = 15.74111793583134
$ws.Range("F14").Value = 40.37256396567957
$ws.Range("G14").Value = 3.667459800632158
$ws.Range("J14").Value = 8.488998814840375
$ws.Range("K14").Value = 10.97908969446897
$ws.Range("L14").Value = 12.3211635536195
$ws.Range("N14").Value = 20.17326490886794
$ws.Range("O14").Value = 26.56139797488807

$ws.Range("B15").Value = 15.68564778269132
$ws.Range("C15").Value = 9.966849439647852
$ws.Range("E15").Value = 15.7385976446824
$ws.Range("F15").Value = 40.36769205178405
$ws.Range("G15").Value = 3.667647274609276
$ws.Range("J15").Value = 8.488755737635611
$ws.Range("K15").Value = 10.96021925614576
$ws.Range("L15").Value = 12.31711319829204
$ws.Range("N15").Value = 20.18004766844074
$ws.Range("O15").Value = 26.5683142208995

$ws.Range("B16").Value = 15.52716461332451
$ws.Range("C16").Value = 9.975186734759275
$ws.Range("E16").Value = 15.72484837219905
$ws.Range("F16").Value = 40.34233172465041
$ws.Range("G16").Value = 3.668738318464456
$ws.Range("J16").Value = 8.487386864443696
$ws.Range("K16").Value = 10.85188891497482
$ws.Range("L16").Value = 12.29438458554826
$ws.Range("N16").Value = 20.21942941723374
$ws.Range("O16").Value = 26.60940147631961

$ws.Range("B17").Value = 15.42979619916395
$ws.Range("C17").Value = 9.980434222568348
$ws.Range("E17").Value = 15.71703519158456
$ws.Range("F17").Value = 40.3290622999107
$ws.Range("G17").Value = 3.669422578556557
$ws.Range("J17").Value = 8.486568716741893
$ws.Range("K17").Value = 10.78529238713256
$ws.Range("L17").Value = 12.28087570637625
$ws.Range("N17").Value = 20.24404710606831
$ws.Range("O17").Value = 26.63590271944162

$ws.Range("B18").Value = 15.37375866022409
$ws.Range("C18").Value = 9.983501353195621
$ws.Range("E18").Value = 15.71277072772457
$ws.Range("F18").Value = 40.32227295174858
$ws.Range("G18").Value = 3.669821650276847
$ws.Range("J18").Value = 8.486106122765646
$ws.Range("K18").Value = 10.74694923431143
$ws.Range("L18").Value = 12.27326694488777
$ws.Range("N18").Value = 20.25837533780856
$ws.Range("O18").Value = 26.651621222048

$ws.Range("B19").Value = 15.35478134295153
$ws.Range("C19").Value = 9.984548241625642
$ws.Range("E19").Value = 15.71136636230693
$ws.Range("F19").Value = 40.32011912987164
$ws.Range("G19").Value = 3.669957715557859
$ws.Range("J19").Value = 8.485950871319179
$ws.Range("K19").Value = 10.73396148362837
$ws.Range("L19").Value = 12.2707185776145
$ws.Range("N19").Value = 20.26325565647414
$ws.Range("O19").Value = 26.6570248964483

$ws.Range("B20").Value = 15.44016516685101
$ws.Range("C20").Value = 9.979870558356085
$ws.Range("E20").Value = 15.71784319184293
$ws.Range("F20").Value = 40.33038765979505
$ws.Range("G20").Value = 3.669349168612742
$ws.Range("J20").Value = 8.486654984671215
$ws.Range("K20").Value = 10.79238598219208
$ws.Range("L20").Value = 12.28229710172156
$ws.Range("N20").Value = 20.24140905066698
$ws.Range("O20").Value = 26.63303237530275

$ws.Range("B21").Value = 15.72650502093922
$ws.Range("C21").Value = 9.964737897026515
$ws.Range("E21").Value = 15.74233871199124
$ws.Range("F21").Value = 40.37494536416769
$ws.Range("G21").Value = 3.667370196828252
$ws.Range("J21").Value = 8.489115801888888
$ws.Range("K21").Value = 10.98813475659357
$ws.Range("L21").Value = 12.32311428417598
$ws.Range("N21").Value = 20.17002144364247
$ws.Range("O21").Value = 26.5581071645923

$ws.Range("B22").Value = 15.91294090217656
$ws.Range("C22").Value = 9.955281271849222
$ws.Range("E22").Value = 15.76036839850507
$ws.Range("F22").Value = 40.41152857246511
$ws.Range("G22").Value = 3.666125556184286
$ws.Range("J22").Value = 8.490794306412912
$ws.Range("K22").Value = 11.11545885984953
$ws.Range("L22").Value = 12.35119247406961
$ws.Range("N22").Value = 20.12486053668769
$ws.Range("O22").Value = 26.51338390402859

$ws.Range("B23").Value = 15.81352008591743
$ws.Range("C23").Value = 9.960288893099065
$ws.Range("E23").Value = 15.75056190541936
$ws.Range("F23").Value = 40.39132295631656
$ws.Range("G23").Value = 3.666785398620193
$ws.Range("J23").Value = 8.489892074352481
$ws.Range("K23").Value = 11.04757191861853
$ws.Range("L23").Value = 12.33608002307412
$ws.Range("N23").Value = 20.14882735414396
$ws.Range("O23").Value = 26.53686474971998

$ws.Range("B24").Value = 15.43547753621847
$ws.Range("C24").Value = 9.980125234310218
$ws.Range("E24").Value = 15.71747718624151
$ws.Range("F24").Value = 40.32978584940314
$ws.Range("G24").Value = 3.669382339548057
$ws.Range("J24").Value = 8.486615958725761
$ws.Range("K24").Value = 10.78917913904966
$ws.Range("L24").Value = 12.28165399738943
$ws.Range("N24").Value = 20.2426011696909
$ws.Range("O24").Value = 26.63432855474856

$ws.Range("B25").Value = 15.02665824816331
$ws.Range("C25").Value = 10.0033933820091
$ws.Range("E25").Value = 15.69057393303463
$ws.Range("F25").Value = 40.29564725125706
$ws.Range("G25").Value = 3.672394140650407
$ws.Range("J25").Value = 8.483387823713857
$ws.Range("K25").Value = 10.50915220426303
$ws.Range("L25").Value = 12.22912639090664
$ws.Range("N25").Value = 20.35020964861654
$ws.Range("O25").Value = 26.75763981222441
